# Re-ran the averaged-intensities notebook including the new spiral
# sampling schemes (Spiral-90deg-10rot-5space, Spiral-90deg-15rot-5space,
# Spiral-90deg-10rot-3space). The spreadsheet table (HKL index / scheme
# name / per-reflection averaged intensities) grows from 14 to 17 schemes,
# and the previously-last "Gaussian-Quadrature" row's recomputed values
# move up to directly follow the "Ring Perpendicular to *" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17:19 do not exist yet in the sheet - grow the table by copying the
# bold/bordered/centered number style used throughout column A (row 16 is
# the current last row) down onto the three new rows before writing values.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9898537016774713
$ws.Range("D10").Value = 0.9479731764878288
$ws.Range("E10").Value = 1.018040081265707
$ws.Range("F10").Value = 0.9898537016774713
$ws.Range("G10").Value = 0.9508196235837006
$ws.Range("H10").Value = 1.086662413002611
$ws.Range("I10").Value = 1.010712039757009
$ws.Range("J10").Value = 0.9479731764878288
$ws.Range("K10").Value = 0.9830066288767679
$ws.Range("L10").Value = 0.9864301652771196
$ws.Range("M10").Value = 1.000676839295721

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9461476045655763
$ws.Range("D11").Value = 0.9497201514044133
$ws.Range("E11").Value = 1.04740182780279
$ws.Range("F11").Value = 0.9461476045655763
$ws.Range("G11").Value = 0.9429022807811757
$ws.Range("H11").Value = 1.171857161548419
$ws.Range("I11").Value = 1.016186540646336
$ws.Range("J11").Value = 0.9497201514044133
$ws.Range("K11").Value = 0.9985609896036018
$ws.Range("L11").Value = 0.9723542970845891
$ws.Range("M11").Value = 1.012369261124785

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.947160214983805
$ws.Range("D12").Value = 0.9453326949076158
$ws.Range("E12").Value = 1.048520564252864
$ws.Range("F12").Value = 0.947160214983805
$ws.Range("G12").Value = 0.9405444506086148
$ws.Range("H12").Value = 1.174589250976436
$ws.Range("I12").Value = 1.017257384990399
$ws.Range("J12").Value = 0.9453326949076158
$ws.Range("K12").Value = 0.9969266295802399
$ws.Range("L12").Value = 0.9720434222820225
$ws.Range("M12").Value = 1.012234093453289

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9462572423490823
$ws.Range("D13").Value = 0.9487236174581086
$ws.Range("E13").Value = 1.047744034250202
$ws.Range("F13").Value = 0.9462572423490823
$ws.Range("G13").Value = 0.9423437630885232
$ws.Range("H13").Value = 1.172695086564939
$ws.Range("I13").Value = 1.016474402639649
$ws.Range("J13").Value = 0.9487236174581086
$ws.Range("K13").Value = 0.9982338258541553
$ws.Range("L13").Value = 0.9722455341016187
$ws.Range("M13").Value = 1.012373024391751

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.225056000000002
$ws.Range("D14").Value = 0.4844079999999993
$ws.Range("E14").Value = 1.079244000000003
$ws.Range("F14").Value = 1.225056000000002
$ws.Range("G14").Value = 0.7270519999999996
$ws.Range("H14").Value = 1.226272000000001
$ws.Range("I14").Value = 1.105956
$ws.Range("J14").Value = 0.4844079999999993
$ws.Range("K14").Value = 0.781826000000001
$ws.Range("L14").Value = 1.003441000000002
$ws.Range("M14").Value = 0.9746646666666674

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.36
$ws.Range("D15").Value = 0.53
$ws.Range("E15").Value = 0.99
$ws.Range("F15").Value = 1.36
$ws.Range("G15").Value = 0.7890249999999992
$ws.Range("H15").Value = 0.96
$ws.Range("I15").Value = 1.07
$ws.Range("J15").Value = 0.53
$ws.Range("K15").Value = 0.76
$ws.Range("L15").Value = 1.06
$ws.Range("M15").Value = 0.9498374999999998

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.204750386380798
$ws.Range("D16").Value = 0.7241841140736001
$ws.Range("E16").Value = 0.9929773148159994
$ws.Range("F16").Value = 1.204750386380798
$ws.Range("G16").Value = 0.8754622837759968
$ws.Range("H16").Value = 0.9778392190976002
$ws.Range("I16").Value = 1.039254409420801
$ws.Range("J16").Value = 0.7241841140736001
$ws.Range("K16").Value = 0.8585807144447997
$ws.Range("L16").Value = 1.031665550412799
$ws.Range("M16").Value = 0.9690779545941325

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9930936015050272
$ws.Range("D17").Value = 0.9955011438546748
$ws.Range("E17").Value = 0.9944649316965294
$ws.Range("F17").Value = 0.9930936015050272
$ws.Range("G17").Value = 0.9938638224076486
$ws.Range("H17").Value = 0.9961850192386988
$ws.Range("I17").Value = 0.9944806121213408
$ws.Range("J17").Value = 0.9955011438546748
$ws.Range("K17").Value = 0.9949830377756022
$ws.Range("L17").Value = 0.9940383196403146
$ws.Range("M17").Value = 0.9945981884706533

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.006919063885568
$ws.Range("D18").Value = 0.9950452198478871
$ws.Range("E18").Value = 0.9869543391756601
$ws.Range("F18").Value = 1.006919063885568
$ws.Range("G18").Value = 1.000119122551058
$ws.Range("H18").Value = 0.9684168653373026
$ws.Range("I18").Value = 0.9932099826403843
$ws.Range("J18").Value = 0.9950452198478871
$ws.Range("K18").Value = 0.9909997795117735
$ws.Range("L18").Value = 0.9989594216986708
$ws.Range("M18").Value = 0.9917774322396434

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9851121695331951
$ws.Range("D19").Value = 1.035032075919488
$ws.Range("E19").Value = 0.9850720154504773
$ws.Range("F19").Value = 0.9851121695331951
$ws.Range("G19").Value = 1.016090513674558
$ws.Range("H19").Value = 0.9732386689609874
$ws.Range("I19").Value = 0.9842497377476701
$ws.Range("J19").Value = 1.035032075919488
$ws.Range("K19").Value = 1.010052045684982
$ws.Range("L19").Value = 0.9975821076090888
$ws.Range("M19").Value = 0.9964658635477294
